$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Data")

# Update formulas to divide by 10 (values were in different units)
$ws.Range("B88").Formula = "=B54/10"
$ws.Range("C88").Formula = "=B55/10"
$ws.Range("B96").Formula = "=B87/10"

# Apply the new number format (2 decimal places) with the same green fill
$ws.Range("B88:C88").NumberFormat = """$""#,##0.00"
$ws.Range("B96").NumberFormat = """$""#,##0.00"
